# quarterly.xlsx update — roll the quarterly columns forward by one
# quarter (E..N now ends on the quarter ending 1401/12 instead of
# 1401/09) and refresh the read_price-derived expense/personnel figures
# for every data row with the newly shifted-in quarter's numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# ---------------------------------------------------------------------
# 1. Column headers (row 8 and row 24) - each quarter label advances by
#    one quarter; a new trailing quarter is introduced in column N.
# ---------------------------------------------------------------------
$headerCols = @("E","F","G","H","I","J","K","L","M","N")
$headers = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "8").Value = $headers[$i]
    $ws.Range($headerCols[$i] + "24").Value = $headers[$i]
}

# ---------------------------------------------------------------------
# 2. Data rows - values shift one column to the left (E<-F, F<-G, ...)
#    and the newly vacated column N receives the freshly read figure.
# ---------------------------------------------------------------------
$dataCols = @("E","F","G","H","I","J","K","L","M","N")
$rowValues = @{
    10 = @(2951, 12776, 3496, 4661, 2240, 5069, 3704, 3750, 2050, 4446)
    13 = @(84, 4122, 0, 1271, 211, 4080, 122, 301, 1154, 6467)
    15 = @(0, 353, 0, 0, 276, 102, 155, 171, 174, 162)
    16 = @(480, 782, 440, 789, 608, 571, 591, 600, 582, 582)
    17 = @(20622, 10674, 19385, 24595, 20340, 29947, 34228, 27412, 32102, 42279)
    19 = @(4245, 14400, 8681, 7984, 8873, 20228, 7153, 10633, 13809, 23257)
    20 = @(28382, 43107, 32002, 39300, 32548, 59997, 45953, 42867, 49871, 77193)
    26 = @(123, 124, 124, 129, 125, 126, 123, 122, 121, 123)
    27 = @(62, 62, 62, 61, 62, 62, 61, 62, 62, 60)
}

foreach ($r in $rowValues.Keys) {
    $vals = $rowValues[$r]
    for ($i = 0; $i -lt $dataCols.Length; $i++) {
        $ws.Range($dataCols[$i] + $r).Value = $vals[$i]
    }
}
